# Applies "Update latest output (run 90)" changes to optimisation_result.xlsx
# Sheet "Schedule": rows 2-4 (schedule summary recompute)
# Sheet "Detailed": rows 45-97 (price forecast -> historical revisions, pump status)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" updates (rows 2-4) ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("B2").Value = 46040.91666666666
$schedule.Range("C2").Value = 15.5
$schedule.Range("D2").Value = 58.59
$schedule.Range("E2").Value = 232.031553
$schedule.Range("F2").Value = 3.960258627752177
$schedule.Range("A3").Value = 46040.97916666666
$schedule.Range("C3").Value = 4
$schedule.Range("D3").Value = 15.12
$schedule.Range("E3").Value = 367.00941225
$schedule.Range("F3").Value = 24.27310927579365
$schedule.Range("E4").Value = -31.19719199999999
$schedule.Range("F4").Value = -0.9709676937441643

# --- Sheet "Detailed" updates (rows 45-97) ---
$detailed = $wb.Worksheets.Item("Detailed")
$detailed.Range("B45").Value = 36.2
$detailed.Range("E45").Value = "ON"
$detailed.Range("B46").Value = 54.93997
$detailed.Range("C47").Value = "historical"
$detailed.Range("B48").Value = 56.98
$detailed.Range("C48").Value = "historical"
$detailed.Range("E48").Value = "OFF"
$detailed.Range("B49").Value = 56.98
$detailed.Range("C49").Value = "historical"
$detailed.Range("B50").Value = 56.98
$detailed.Range("B51").Value = 56.98
$detailed.Range("B52").Value = 56.98
$detailed.Range("B53").Value = 40.54
$detailed.Range("B54").Value = 35.87995
$detailed.Range("B55").Value = 35.87996
$detailed.Range("B56").Value = 36.2
$detailed.Range("B57").Value = 48.78448
$detailed.Range("B59").Value = 58.3247
$detailed.Range("B60").Value = 58.23997
$detailed.Range("B61").Value = 60.48315
$detailed.Range("B64").Value = 35.88
$detailed.Range("B65").Value = 4.32936
$detailed.Range("B66").Value = 0.7
$detailed.Range("B67").Value = 0.7
$detailed.Range("B68").Value = 0.51
$detailed.Range("B69").Value = -5.01
$detailed.Range("B70").Value = -6.05265
$detailed.Range("B71").Value = -5.95647
$detailed.Range("B72").Value = -0.93185
$detailed.Range("B73").Value = -0.92598
$detailed.Range("B74").Value = -4.65702
$detailed.Range("B75").Value = -4.67425
$detailed.Range("B76").Value = -5.01
$detailed.Range("B77").Value = -4.78417
$detailed.Range("B78").Value = -0.9435
$detailed.Range("B79").Value = 0.00002
$detailed.Range("B80").Value = 0.009390000000000001
$detailed.Range("B81").Value = 0.7
$detailed.Range("B82").Value = 0.7
$detailed.Range("B83").Value = -2.47065
$detailed.Range("B84").Value = -0.53799
$detailed.Range("B85").Value = 0.33937
$detailed.Range("B86").Value = 11.98812
$detailed.Range("B87").Value = 51.96257
$detailed.Range("B88").Value = 56.69797
$detailed.Range("B89").Value = 64.02478000000001
$detailed.Range("B90").Value = 64.08503
$detailed.Range("B91").Value = 65
$detailed.Range("B94").Value = 61.4952
$detailed.Range("B95").Value = 59.32479
$detailed.Range("B96").Value = 58.48496
$detailed.Range("B97").Value = 61.96632
